# Applies scheduled-runner market price/profit updates to the Leve profit sheets.
# Generated from the authoritative cell-level diff between before/after workbook states.
$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 765.4
$ws.Range("I6").Value = 765.4
$ws.Range("K6").Value = 2296.2
$ws.Range("M6").Value = -2184.2

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 30069.5
$ws.Range("I11").Value = 30069.5
$ws.Range("K11").Value = 30069.5
$ws.Range("M11").Value = -29929.5

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 45466.305
$ws.Range("J17").Value = 47487.5
$ws.Range("L17").Value = 142462.5
$ws.Range("N17").Value = -142798.5

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 767.0952
$ws.Range("I42").Value = 100.64286
$ws.Range("J42").Value = 2100
$ws.Range("K42").Value = 301.92858
$ws.Range("L42").Value = 6300
$ws.Range("M42").Value = -71.92858000000001
$ws.Range("N42").Value = -6760

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4991.3335
$ws.Range("I62").Value = 4996
$ws.Range("K62").Value = 4996
$ws.Range("M62").Value = -4372

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4991.3335
$ws.Range("I65").Value = 4996
$ws.Range("K65").Value = 24980
$ws.Range("M65").Value = -21860

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3490.9092
$ws.Range("I88").Value = 880
$ws.Range("J88").Value = 5666.6665
$ws.Range("K88").Value = 880
$ws.Range("L88").Value = 5666.6665
$ws.Range("M88").Value = -474
$ws.Range("N88").Value = -6478.6665

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3490.9092
$ws.Range("I91").Value = 880
$ws.Range("J91").Value = 5666.6665
$ws.Range("K91").Value = 880
$ws.Range("L91").Value = 5666.6665
$ws.Range("M91").Value = 524
$ws.Range("N91").Value = -8474.666499999999

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1663.75
$ws.Range("I96").Value = 1292.7273
$ws.Range("J96").Value = 2480
$ws.Range("K96").Value = 3878.1819
$ws.Range("L96").Value = 7440
$ws.Range("M96").Value = -2505.1819
$ws.Range("N96").Value = -10186

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1397.3334
$ws.Range("I100").Value = 1397.3334
$ws.Range("K100").Value = 1397.3334
$ws.Range("M100").Value = -856.3334

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1295
$ws.Range("I111").Value = 1295
$ws.Range("K111").Value = 3885
$ws.Range("M111").Value = -818

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4666.6665
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 849961.9
$ws.Range("I61").Value = 23054.818
$ws.Range("J61").Value = 3525249.2
$ws.Range("K61").Value = 23054.818
$ws.Range("L61").Value = 3525249.2
$ws.Range("M61").Value = -22842.818
$ws.Range("N61").Value = -3525673.2

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 316055.6
$ws.Range("I74").Value = 2051.1968
$ws.Range("K74").Value = 2051.1968
$ws.Range("M74").Value = -1177.1968

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 316055.6
$ws.Range("I77").Value = 2051.1968
$ws.Range("K77").Value = 10255.984
$ws.Range("M77").Value = -5887.984

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1891.1428
$ws.Range("I110").Value = 1891.1428
$ws.Range("K110").Value = 1891.1428
$ws.Range("M110").Value = 153.8571999999999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3651.611
$ws.Range("I132").Value = 2297.25
$ws.Range("K132").Value = 6891.75
$ws.Range("M132").Value = -4361.75

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 849961.9
$ws.Range("I136").Value = 23054.818
$ws.Range("J136").Value = 3525249.2
$ws.Range("K136").Value = 69164.454
$ws.Range("L136").Value = 10575747.6
$ws.Range("M136").Value = -66614.454
$ws.Range("N136").Value = -10580847.6

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7492.933
$ws.Range("I99").Value = 7849.5713
$ws.Range("K99").Value = 7849.5713
$ws.Range("M99").Value = -6351.5713

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9448.212
$ws.Range("I105").Value = 7368.0454
$ws.Range("K105").Value = 7368.0454
$ws.Range("M105").Value = -5621.0454

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6891.1143
$ws.Range("I107").Value = 7126.1816
$ws.Range("K107").Value = 7126.1816
$ws.Range("M107").Value = -5206.1816

# BSM row 128
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 10000
$ws.Range("I128").Value = 10000
$ws.Range("K128").Value = 30000
$ws.Range("M128").Value = -27510

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33335142
$ws.Range("I134").Value = 1701.409
$ws.Range("K134").Value = 5104.227000000001
$ws.Range("M134").Value = -2569.227000000001

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1004
$ws.Range("I2").Value = 1004
$ws.Range("K2").Value = 1004
$ws.Range("M2").Value = -891

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 592.8387
$ws.Range("I22").Value = 567.12
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 567.12
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -217.12
$ws.Range("N22").Value = -1400

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2830
$ws.Range("I132").Value = 1947.5
$ws.Range("J132").Value = 3334.2856
$ws.Range("K132").Value = 5842.5
$ws.Range("L132").Value = 10002.8568
$ws.Range("M132").Value = -3312.5
$ws.Range("N132").Value = -15062.8568

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 786533.2
$ws.Range("I4").Value = 866.75
$ws.Range("J4").Value = 1100799.8
$ws.Range("K4").Value = 2600.25
$ws.Range("L4").Value = 3302399.4
$ws.Range("M4").Value = -2488.25
$ws.Range("N4").Value = -3302623.4

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2676982
$ws.Range("I131").Value = 18183918
$ws.Range("J131").Value = 3372.4138
$ws.Range("K131").Value = 54551754
$ws.Range("L131").Value = 10117.2414
$ws.Range("M131").Value = -54546714
$ws.Range("N131").Value = -20197.2414

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12633375
$ws.Range("J80").Value = 50055200
$ws.Range("L80").Value = 50055200
$ws.Range("N80").Value = -50057196

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 12633375
$ws.Range("J83").Value = 50055200
$ws.Range("L83").Value = 250276000
$ws.Range("N83").Value = -250285984

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 680.85
$ws.Range("I97").Value = 494.6154
$ws.Range("K97").Value = 494.6154
$ws.Range("M97").Value = 1.384599999999978

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1070163.4
$ws.Range("I132").Value = 1498
$ws.Range("K132").Value = 4494
$ws.Range("M132").Value = -1964

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3578.4285
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3578.4285
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3578.4285
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -5076.4285

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3578.4285
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3578.4285
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 17892.1425
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -25380.1425

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3980.7778
$ws.Range("I100").Value = 3945.5715
$ws.Range("J100").Value = 3993.1
$ws.Range("K100").Value = 3945.5715
$ws.Range("L100").Value = 3993.1
$ws.Range("M100").Value = -3404.5715
$ws.Range("N100").Value = -5075.1

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -88

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 32260210
$ws.Range("I132").Value = 66668228
$ws.Range("J132").Value = 2692.8125
$ws.Range("K132").Value = 200004684
$ws.Range("L132").Value = 8078.4375
$ws.Range("M132").Value = -200002154
$ws.Range("N132").Value = -13138.4375
